$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 1232, 0.000473, 0.000459, 0.02918),
    @(3, 1002, 0.000834, 0.001291, -0.548528),
    @(4, 1094, 0.000223, 0.000244, -0.09286800000000001),
    @(5, 1163, 0.000124, 0.000063, 0.490702),
    @(6, 1252, 0.000002, 0.000005, -1.240692),
    @(7, 1068, 0.000277, 0.000287, -0.03326),
    @(8, 1016, 0.000586, 0.000654, -0.116363),
    @(9, 1052, 0.001344, 0.00114, 0.151605),
    @(10, 1281, 0.000067, 0.000079, -0.174636),
    @(11, 1195, 0.000006, 0.000007, -0.07692400000000001),
    @(12, 1128, 0.000644, 0.000691, -0.073416),
    @(13, 1048, 0.000846, 0.000679, 0.197235),
    @(14, 1064, 0.000167, 0.000224, -0.342878),
    @(15, 1003, 0.000247, 0.000285, -0.152259),
    @(16, 1076, 0.000727, 0.000776, -0.067181),
    @(17, 1071, 0.00001, 0.000008, 0.197509),
    @(18, 1069, 0.000112, 0.000101, 0.099215),
    @(19, 1224, 0.000108, 0.00008899999999999999, 0.175884),
    @(20, 1151, 0.000163, 0.000079, 0.516611),
    @(21, 1102, 0.000662, 0.000725, -0.09388100000000001),
    @(22, 1037, 0.003112, 0.003639, -0.169555),
    @(23, 1228, 0.000318, 0.000201, 0.368786),
    @(24, 1038, 0.001427, 0.00163, -0.142086),
    @(25, 1159, 0.000017, 0.000029, -0.756073),
    @(26, 1261, 0.000005, 0.000009, -0.729257),
    @(27, 1177, 0.000056, 0.000044, 0.214269),
    @(28, 1267, 0.000087, 0.000098, -0.13214),
    @(29, 1087, 0.000011, 0.00001, 0.099707),
    @(30, 1188, 0.000298, 0.000213, 0.284888),
    @(31, 1011, 0.00121, 0.001883, -0.556957),
    @(32, 1107, 0.000032, 0.000038, -0.20312),
    @(33, 1130, 0.000063, 0.00006999999999999999, -0.10626),
    @(34, 1035, 0.000454, 0.000454, 0.000464),
    @(35, 1203, 0.00005, 0.000043, 0.149672),
    @(36, 1115, 0.0004, 0.00061, -0.527344),
    @(37, 1110, 0.0004, 0.0004, 0.000519),
    @(38, 1025, 0.000633, 0.000678, -0.07016500000000001),
    @(39, 1280, 0.000293, 0.000302, -0.032419),
    @(40, 1139, 0.000013, 0.000001, 0.889523),
    @(41, 1260, 0.000025, 0.000052, -1.091652),
    @(42, 1104, 0.000086, 0.000098, -0.132271),
    @(43, 1208, 0.000011, 0.000028, -1.563135),
    @(44, 1189, 0.00008500000000000001, 0.00006600000000000001, 0.229156),
    @(45, 1190, 0.000008, 0.000012, -0.53796),
    @(46, 1162, 0.000341, 0.000229, 0.329014),
    @(47, 1245, 0.000279, 0.000429, -0.5371359999999999),
    @(48, 1043, 0.001168, 0.000919, 0.213468),
    @(49, 1100, 0.000031, 0.000032, -0.013936),
    @(50, 1045, 0.00028, 0.000203, 0.274172),
    @(51, 1047, 0.001671, 0.001325, 0.206901),
    @(52, 1060, 0.000096, 0.000082, 0.146878),
    @(53, 1242, 0.000021, 0.00002, 0.046048),
    @(54, 1194, 0.000032, 0.000052, -0.592856),
    @(55, 1077, 0.000327, 0.000344, -0.052504),
    @(56, 1133, 0.000041, 0.000017, 0.58345),
    @(57, 1277, 0.000017, 0.000026, -0.5570850000000001),
    @(58, 1073, 0.000209, 0.00019, 0.08673699999999999),
    @(59, 1202, 0.000244, 0.000173, 0.287818),
    @(60, 1086, 0.00006999999999999999, 0.000053, 0.249575),
    @(61, 1182, 0.000015, 0.000029, -0.9417990000000001),
    @(62, 1030, 0.000381, 0.0004929999999999999, -0.294847)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = "[]"
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
}
